$wb = $excel.ActiveWorkbook

# --- Sheet1 (TableToDicts): update selection, it will no longer be the active/selected tab ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("N25").Select() | Out-Null

# --- Add new worksheet "GroupOnIndentations" right after TableToDicts ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "GroupOnIndentations"

# Column widths for columns E and F
$ws2.Columns("E:F").ColumnWidth = 8.88671875

# --- Populate row 5 (header-like row with indentation demo data) ---
$ws2.Range("D5").Value = "x"
$ws2.Range("E5").Value = "y"
$ws2.Range("F5").Value = "z"
$ws2.Range("G5").Value = "zz"

# --- Populate column C rows 6-12 ---
$ws2.Range("C6").Value = "a"
$ws2.Range("C7").Value = "b"
$ws2.Range("C8").Value = "c"
$ws2.Range("C9").Value = "d"
$ws2.Range("C10").Value = "e"
$ws2.Range("C11").Value = "f"
$ws2.Range("C12").Value = "g"

# --- Apply indentation level 1 style: C7, C8, C10, E5 ---
$ws2.Range("C7").HorizontalAlignment = -4131
$ws2.Range("C7").IndentLevel = 1
$ws2.Range("C7").Copy()
$ws2.Range("C8").PasteSpecial(-4122)
$ws2.Range("C10").PasteSpecial(-4122)
$ws2.Range("E5").PasteSpecial(-4122)

# --- Apply indentation level 2 style: C11, C12, F5 ---
$ws2.Range("C11").HorizontalAlignment = -4131
$ws2.Range("C11").IndentLevel = 2
$ws2.Range("C11").Copy()
$ws2.Range("C12").PasteSpecial(-4122)
$ws2.Range("F5").PasteSpecial(-4122)

# --- Make the new sheet the active one and set its selection ---
$ws2.Activate() | Out-Null
$ws2.Range("R21").Select() | Out-Null

# --- Defined names used by tests ---
$wb.Names.Add("__TestGroupColumnsOnIndentations__", "=GroupOnIndentations!`$D`$5:`$G`$5")
$wb.Names.Add("__TestGroupRowsOnIndentations__", "=GroupOnIndentations!`$C`$6:`$C`$12")
